$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 1592.5
$ws.Range("I12").Value = 1906.8334
$ws.Range("J12").Value = 649.5
$ws.Range("K12").Value = 1906.8334
$ws.Range("L12").Value = 649.5
$ws.Range("M12").Value = -1736.8334
$ws.Range("N12").Value = -989.5
# Row 86
$ws.Range("H86").Value = 7697.1177
$ws.Range("I86").Value = 1840.625
$ws.Range("J86").Value = 12902.889
$ws.Range("K86").Value = 1840.625
$ws.Range("L86").Value = 12902.889
$ws.Range("M86").Value = -717.625
$ws.Range("N86").Value = -15148.889
# Row 89
$ws.Range("H89").Value = 7697.1177
$ws.Range("I89").Value = 1840.625
$ws.Range("J89").Value = 12902.889
$ws.Range("K89").Value = 9203.125
$ws.Range("L89").Value = 64514.44499999999
$ws.Range("M89").Value = -3587.125
$ws.Range("N89").Value = -75746.44499999999
# Row 98
$ws.Range("H98").Value = 373.33334
$ws.Range("I98").Value = 295
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 295
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = 1203
$ws.Range("N98").Value = -3996
# Row 116
$ws.Range("H116").Value = 5206.643
$ws.Range("I116").Value = 2998.1428
$ws.Range("J116").Value = 7415.143
$ws.Range("K116").Value = 2998.1428
$ws.Range("L116").Value = 7415.143
$ws.Range("M116").Value = 443.8571999999999
# Row 122
$ws.Range("H122").Value = 373.33334
$ws.Range("I122").Value = 295
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 885
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = 1565
$ws.Range("N122").Value = -7900
# Row 125
$ws.Range("H125").Value = 2003.375
$ws.Range("I125").Value = 1712.8
$ws.Range("J125").Value = 2487.6667
$ws.Range("K125").Value = 15415.2
$ws.Range("L125").Value = 22389.0003
$ws.Range("M125").Value = -12955.2
$ws.Range("N125").Value = -27309.0003
# Row 129
$ws.Range("H129").Value = 1130.3805
$ws.Range("I129").Value = 597
$ws.Range("J129").Value = 1148.3595
$ws.Range("K129").Value = 1791
$ws.Range("L129").Value = 3445.0785
$ws.Range("M129").Value = 3209
$ws.Range("N129").Value = -13445.0785
# Row 137
$ws.Range("H137").Value = 37315.105
$ws.Range("I137").Value = 1374.6842
$ws.Range("J137").Value = 113189.336
$ws.Range("K137").Value = 4124.0526
$ws.Range("L137").Value = 339568.008
$ws.Range("M137").Value = -1574.0526
$ws.Range("N137").Value = -344668.008

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3399.8838
$ws.Range("I45").Value = 2890.1365
$ws.Range("J45").Value = 3933.9048
$ws.Range("K45").Value = 2890.1365
$ws.Range("L45").Value = 3933.9048
$ws.Range("M45").Value = -2513.1365
$ws.Range("N45").Value = -4687.9048
# Row 122
$ws.Range("H122").Value = 2371
$ws.Range("I122").Value = 2364.55
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 7093.650000000001
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -4643.650000000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -150
$ws.Range("N22").ClearContents()
# Row 68
$ws.Range("H68").Value = 64407.777
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 64407.777
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 64407.777
$ws.Range("N68").Value = -65905.777
# Row 71
$ws.Range("H71").Value = 64407.777
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 64407.777
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 193223.331
$ws.Range("N71").Value = -200711.331
# Row 99
$ws.Range("H99").Value = 20119150
$ws.Range("I99").Value = 5211599
$ws.Range("J99").Value = 38466908
$ws.Range("K99").Value = 5211599
$ws.Range("L99").Value = 38466908
$ws.Range("M99").Value = -5210101
# Row 122
$ws.Range("H122").Value = 1086.7646
$ws.Range("I122").Value = 923.8182
$ws.Range("J122").Value = 1385.5
$ws.Range("K122").Value = 2771.4546
$ws.Range("L122").Value = 4156.5
$ws.Range("M122").Value = -321.4546
$ws.Range("N122").Value = -9056.5
# Row 126
$ws.Range("H126").Value = 20119150
$ws.Range("I126").Value = 5211599
$ws.Range("J126").Value = 38466908
$ws.Range("K126").Value = 15634797
$ws.Range("L126").Value = 115400724
$ws.Range("M126").Value = -15632327

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 2067.6123
$ws.Range("I2").Value = 2464.1707
$ws.Range("J2").Value = 35.25
$ws.Range("K2").Value = 14785.0242
$ws.Range("L2").Value = 211.5
$ws.Range("M2").Value = -14672.0242
# Row 131
$ws.Range("H131").Value = 791.38
$ws.Range("I131").Value = 577
$ws.Range("J131").Value = 802.66315
$ws.Range("K131").Value = 1731
$ws.Range("L131").Value = 2407.98945
$ws.Range("M131").Value = 3309
$ws.Range("N131").Value = -12487.98945
# Row 138
$ws.Range("H138").Value = 126315.664
$ws.Range("I138").Value = 1305.4
$ws.Range("J138").Value = 334666.12
$ws.Range("K138").Value = 3916.2
$ws.Range("L138").Value = 1003998.36
$ws.Range("M138").Value = 1223.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 22729418
$ws.Range("I102").Value = 27779954
$ws.Range("J102").Value = 2003.25
$ws.Range("K102").Value = 27779954
$ws.Range("L102").Value = 2003.25
$ws.Range("M102").Value = -27778332
$ws.Range("N102").Value = -5247.25
# Row 122
$ws.Range("H122").Value = 2188
$ws.Range("I122").Value = 1850.875
$ws.Range("J122").Value = 2412.75
$ws.Range("K122").Value = 5552.625
$ws.Range("L122").Value = 7238.25
$ws.Range("M122").Value = -3102.625
# Row 132
$ws.Range("H132").Value = 68169
$ws.Range("I132").Value = 55677.58
$ws.Range("J132").Value = 127503.25
$ws.Range("K132").Value = 167032.74
$ws.Range("L132").Value = 382509.75
$ws.Range("M132").Value = -164502.74

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5431.4375
$ws.Range("I7").Value = 5443.0713
$ws.Range("J7").Value = 5350
$ws.Range("K7").Value = 5443.0713
$ws.Range("L7").Value = 5350
$ws.Range("M7").Value = -5331.0713
# Row 40
$ws.Range("H40").Value = 3515.2354
$ws.Range("I40").Value = 2725.6667
$ws.Range("J40").Value = 3945.9092
$ws.Range("K40").Value = 2725.6667
$ws.Range("L40").Value = 3945.9092
$ws.Range("M40").Value = -2589.6667
$ws.Range("N40").Value = -4217.9092
# Row 55
$ws.Range("H55").Value = 174.3
$ws.Range("I55").Value = 44.2
$ws.Range("J55").Value = 304.4
$ws.Range("K55").Value = 44.2
$ws.Range("L55").Value = 304.4
$ws.Range("M55").Value = 128.8
$ws.Range("N55").Value = -650.4
# Row 70
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
# Row 122
$ws.Range("H122").Value = 1964004.6
$ws.Range("I122").Value = 2181371.8
$ws.Range("J122").Value = 7700
$ws.Range("K122").Value = 6544115.399999999
$ws.Range("L122").Value = 23100
$ws.Range("M122").Value = -6541665.399999999
# Row 126
$ws.Range("H126").Value = 5431.4375
$ws.Range("I126").Value = 5443.0713
$ws.Range("J126").Value = 5350
$ws.Range("K126").Value = 16329.2139
$ws.Range("L126").Value = 16050
$ws.Range("M126").Value = -13859.2139
# Row 136
$ws.Range("H136").Value = 42939.832
$ws.Range("I136").Value = 46572.547
$ws.Range("J136").Value = 2980
$ws.Range("K136").Value = 139717.641
$ws.Range("L136").Value = 8940
$ws.Range("M136").Value = -137167.641

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1942.0476
$ws.Range("I122").Value = 1944.1578
$ws.Range("J122").Value = 1922
$ws.Range("K122").Value = 5832.4734
$ws.Range("L122").Value = 5766
$ws.Range("M122").Value = -3382.4734
# Row 126
$ws.Range("H126").Value = 1543.8889
$ws.Range("I126").Value = 1100
$ws.Range("J126").Value = 1670.7142
$ws.Range("K126").Value = 3300
$ws.Range("L126").Value = 5012.142599999999
$ws.Range("M126").Value = -830
$ws.Range("N126").Value = -9952.142599999999
